$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Month labels: capitalize first letter (jan -> Jan, feb -> Feb, ...)
$months = @("Jan", "Feb", "Mar", "Apr", "May", "Jun", "Jul", "Aug", "Sep", "Oct", "Nov", "Dec")

# Mean / SEM / Std values rounded to 4 decimal places
$meanVals = @(1766.8202, 1661.4328, 1392.2464, 1684.1108, 1567.7157, 1212.0934, 1434.2319, 1254.7059, 1366.844, 1629.1399, 1510.7433, 1476.4311)
$semVals  = @(229.2165, 284.9497, 195.7472, 210.2086, 261.7806, 152.3768, 204.7009, 141.289, 172.4449, 209.6539, 215.2498, 164.9748)
$stdVals  = @(2150.2412, 2596.0155, 1986.6167, 2153.9974, 1775.4822, 834.6023, 1518.1028, 1076.0252, 1335.7522, 1677.2315, 1653.3647, 1466.3281)

for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $months[$i]
    $ws.Cells.Item($row, 4).Value = $meanVals[$i]
    $ws.Cells.Item($row, 5).Value = $semVals[$i]
    $ws.Cells.Item($row, 6).Value = $stdVals[$i]
}
